$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.559.30'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '3.033.56'
$ws.Range('E3').Value = '  -4.41%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.25%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.030.44'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.499'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.83%  '
$ws.Range('E10').Value = '  -3.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.24'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.434'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.44%  '
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.37'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').Value = '3.536.33'
$ws.Range('E16').Value = '  -4.32%  '
$ws.Range('D17').Value = '61.664.07'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '3.031.77'
$ws.Range('E18').Value = '  -4.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '442.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.34%  '
$ws.Range('E22').Value = '  -4.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.14%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.40'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.37'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0957'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.78%  '
$ws.Range('E34').Value = '  -2.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.969'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.65'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '50.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.46%  '
$ws.Range('D38').Value = '0.0₃0694'
$ws.Range('E38').Value = '  -0.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0370'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.85'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('E41').Value = '  -2.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.52'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '376.17'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.00%  '
$ws.Range('D44').Value = '2.670.79'
$ws.Range('E44').Value = '  -4.41%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '123.23'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.237'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '33.83'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.99'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.107'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.66'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.29%  '
